# Update LR-pair data with new TPM-based values.
# The "Sending cluster" / "Target cluster" combinations are expanded from
# 4 rows (all sent from MuSCs) to 6 rows (sent from ECs and from MuSCs,
# each against targets ECs / FAPs / MuSCs), and every numeric metric is
# refreshed to match the new TPM computation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T the numeric metrics (Ligand-expressing cells ... Edge total expression derived specificity)
$rows = @(
    @{ Row = 2;  A = "ECs";   B = "Fgf16"; C = "Fgfr2"; D = "ECs";
       E = 2; F = 1; G = 0.2355055; H = 0.471011; I = 0.2095102068231121; J = 0.2095102068231121;
       K = 1; L = 0.5; M = 0.09207700000000001; N = 0.184154; O = 0.0789959771480734; P = 0.05545240531440215;
       Q = 0.0216846399235; R = 0.086738559694; S = 0.0165504635104867; T = 0.01161784490625944 },

    @{ Row = 3;  A = "ECs";   B = "Fgf16"; C = "Fgfr2"; D = "FAPs";
       E = 2; F = 1; G = 0.2355055; H = 0.471011; I = 0.2095102068231121; J = 0.2095102068231121;
       K = 3; L = 1; M = 0.9897559999999999; N = 2.969268; O = 0.8491451975864605; P = 0.8941052196698643;
       Q = 0.233092981658; R = 1.398557889948; S = 0.1779045859691918; T = 0.1873241694946574 },

    @{ Row = 4;  A = "ECs";   B = "Fgf16"; C = "Fgfr2"; D = "MuSCs";
       E = 2; F = 1; G = 0.2355055; H = 0.471011; I = 0.2095102068231121; J = 0.2095102068231121;
       K = 2; L = 1; M = 0.083758; N = 0.167516; O = 0.07185882526546619; P = 0.05044237501573352;
       Q = 0.019725469669; R = 0.078901878676; S = 0.0150551573434337; T = 0.01056819242219531 },

    @{ Row = 5;  A = "MuSCs"; B = "Fgf16"; C = "Fgfr2"; D = "ECs";
       E = 2; F = 1; G = 0.888571; H = 1.777142; I = 0.7904897931768879; J = 0.7904897931768879;
       K = 1; L = 0.5; M = 0.09207700000000001; N = 0.184154; O = 0.0789959771480734; P = 0.05545240531440215;
       Q = 0.08181695196700001; R = 0.327267807868; S = 0.06244551363758671; T = 0.04383456040814272 },

    @{ Row = 6;  A = "MuSCs"; B = "Fgf16"; C = "Fgfr2"; D = "FAPs";
       E = 2; F = 1; G = 0.888571; H = 1.777142; I = 0.7904897931768879; J = 0.7904897931768879;
       K = 3; L = 1; M = 0.9897559999999999; N = 2.969268; O = 0.8491451975864605; P = 0.8941052196698643;
       Q = 0.8794684786759999; R = 5.276810872055999; S = 0.6712406116172688; T = 0.706781050175207 },

    @{ Row = 7;  A = "MuSCs"; B = "Fgf16"; C = "Fgfr2"; D = "MuSCs";
       E = 2; F = 1; G = 0.888571; H = 1.777142; I = 0.7904897931768879; J = 0.7904897931768879;
       K = 2; L = 1; M = 0.083758; N = 0.167516; O = 0.07185882526546619; P = 0.05044237501573352;
       Q = 0.074424929818; R = 0.297699719272; S = 0.0568036679220325; T = 0.03987418259353821 }
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in $rows) {
    $rowNum = $r.Row
    foreach ($col in $columns) {
        $ws.Range("$col$rowNum").Value = $r[$col]
    }
}
